# Applies the "Changes to UX folder" commit:
#  - rename worksheet "Sheet1" -> "Graphs"
#  - repoint every chart series (on the renamed sheet) from Sheet1!... to Graphs!...
#  - drop the stale _xlchart.v2.* defined names (Excel's chart-insert leftovers)
#  - make "Graphs" the active / selected tab instead of "Form responses 1"

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Sheet1" worksheet to "Graphs" -------------------------
$graphs = $wb.Worksheets.Item("Sheet1")
$graphs.Name = "Graphs"

# --- 2. Fix up each embedded chart's series formulas so they reference ----
#        the renamed sheet instead of the old "Sheet1" name.
foreach ($co in $graphs.ChartObjects()) {
    $chart = $co.Chart
    foreach ($s in $chart.SeriesCollection()) {
        $s.Formula = $s.Formula -replace "Sheet1!", "Graphs!"
    }
}

# --- 3. Remove the leftover hidden "_xlchart" defined names ---------------
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# --- 4. Make "Graphs" the selected/active sheet instead of the form sheet -
$graphs.Activate()
